# Reorder the names/emails listed in the "Recorded By" column (G) so that
# "System" (or an admin/backup account) no longer sorts first.
#
# Observed substitutions (exact, case-sensitive string matches):
#   "backup@backdoor.com, System, system"  -> "backup@backdoor.com, system, System"
#   "System, dnasr281@gmail.com"           -> "dnasr281@gmail.com, System"
#   "admin@admin.com, dnasr281@gmail.com"  -> "dnasr281@gmail.com, admin@admin.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$map = @{
    "backup@backdoor.com, System, system" = "backup@backdoor.com, system, System";
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System";
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com";
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text
    if ($null -ne $val -and $map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
}
